# Fill in the previously-blank "Titus" evaluation block (rows 13-17,
# columns C:Z) with the scores from the completed human-evaluation
# spreadsheet. Formulas in rows 28-32 (AVERAGE across the five raters)
# recalculate automatically once these cells are populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-ColLetter($col) {
    $letter = ""
    while ($col -gt 0) {
        $rem = ($col - 1) % 26
        $letter = [char](65 + $rem) + $letter
        $col = [int](($col - $rem - 1) / 26)
    }
    return $letter
}

function Set-RowValues($Sheet, $Row, $Values) {
    $n = $Values.Length
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) { $arr[0,$i] = $Values[$i] }
    $startCol = 3   # column C
    $endCol = $startCol + $n - 1
    $addr = (Get-ColLetter $startCol) + $Row + ":" + (Get-ColLetter $endCol) + $Row
    $Sheet.Range($addr).Value = $arr
}

# Columns C:T are per-transcript scores (all 1's), columns U:Z are the
# six "Overall" rating columns for this rater.
Set-RowValues $ws 13 @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1, 3,3,2.5,3,3,4)
Set-RowValues $ws 14 @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1, 3,3,2.5,4,3,3)
Set-RowValues $ws 15 @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1, 3,2.5,2,3,3,3)
Set-RowValues $ws 16 @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1, 2.5,2,2,3,3,3)
Set-RowValues $ws 17 @(1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1, 3,2.5,2.5,3,3,3)

# Match the on-disk selection state left behind by the author: the
# block that was just filled in (C13:T17) ends up selected.
$ws.Range("C13:T17").Select() | Out-Null
